# Updated comments to archive jobs
# Adds a new "Comments" column (R) to Sheet1, flagging most jobs as
# "Archive" and two as "Dee to Review"; also groups/hides the now
# low-priority detail columns (D:M) under an outline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Comments" header (bold, yellow fill) -----------------------------
$ws.Range("R1").Value = "Comments"
$ws.Range("R1").Font.Bold = $true
$ws.Range("R1").Interior.Color = 65535

# --- Per-row comment values -------------------------------------------------
$archiveRows = @(3,4,5,6,7,10,11,12,13,18,19,20,21,22,24,28,29,30,31,45)
foreach ($r in $archiveRows) {
    $ws.Range("R$r").Value = "Archive"
}

$reviewRows = @(23,26)
foreach ($r in $reviewRows) {
    $ws.Range("R$r").Value = "Dee to Review"
}

# --- Group and hide the now-secondary detail columns (D:M) -----------------
$detailCols = $ws.Range("D1:M1").EntireColumn
$detailCols.Group()
$detailCols.Hidden = $true

# --- Restore the active selection/view --------------------------------------
$ws.Range("R46").Select()
